# Update the "想去人数" (F column) counts that were refreshed when the
# gh-pages output was regenerated. The same updates apply to both the
# "展览" sheet (rows keyed by their original row number) and the
# "全部类型" sheet (which merges rows from 展览/演出/本地生活, so its row
# numbers differ slightly).

$wb = $excel.ActiveWorkbook

# Updates for sheet "展览": row number -> new F value
$exhibitionUpdates = @{
    2  = 12838
    3  = 628
    6  = 323
    7  = 403
    9  = 12845
    10 = 39
    11 = 21
    12 = 5236
    18 = 37
    20 = 675
    22 = 6174
    23 = 1156
    24 = 3624
    26 = 44
}

# Updates for sheet "全部类型": row number -> new F value
$allTypesUpdates = @{
    2  = 12838
    3  = 628
    6  = 323
    8  = 403
    10 = 12845
    11 = 39
    12 = 21
    13 = 5236
    19 = 37
    21 = 675
    24 = 6174
    25 = 1156
    26 = 3624
    28 = 44
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
